$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Update existing parameter bounds ---
$ws.Range("C2").Value = 100

$ws.Range("B3").Value = 0.1
$ws.Range("C3").Value = 100

$ws.Range("B4").Value = 0.00001

$ws.Range("B5").Value = 0.00001

$ws.Range("C6").Value = 10

$ws.Range("B8").Value = 0.5
$ws.Range("C8").Value = 20

$ws.Range("C12").Value = 1

$ws.Range("B15").Value = 0.1
$ws.Range("C15").Value = 100

$ws.Range("C17").Value = 0.1

$ws.Range("B19").Value = 0.1
$ws.Range("C19").Value = 100

$ws.Range("B20").Value = 0.1
$ws.Range("C20").Value = 100

$ws.Range("B21").Value = 0.000001

$ws.Range("B22").Value = 0.001

$ws.Range("B23").Value = 0.5
$ws.Range("C23").Value = 20

# --- Add new row 24 (new circuit: LacI_rep_3mut_P3) ---
# Seed formatting by copying the row above (same column layout/number formats),
# then overwrite the contents with the new circuit's values.
$ws.Range("A17:G17").Copy($ws.Range("A24:G24"))

# (G24's shared string is written first so it lands before A24's in sharedStrings.xml,
#  matching the order the workbook was authored in)
$ws.Range("G24").Value = "theta_{LacI_W220F_Q60G_T167A-pt7}"
$ws.Range("A24").Value = "LacI_rep_3mut_P3"
$ws.Range("B24").Value = 0.00001
$ws.Range("C24").Value = 0.1
$ws.Range("D24").Value = 1
$ws.Range("E24").Value = 1
$ws.Range("F24").Value = "yes"

# --- View state: scroll down a bit and select the next empty row ---
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 3
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A25").Select()
